$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")
$lo = $ws.ListObjects.Item("Snippets")

# The "Snippets" table currently spans A1:E60 (header + 59 data rows).
# Insert a brand-new data row right above the final existing row (old
# row 60, the "Table.getCell" snippet) so that row becomes row 61 and
# the new "Style" enum snippet becomes row 60.
$ws.Rows(60).Insert()
$lo.Resize($ws.Range("A1:E61"))

$ws.Range("A60").Value = "Style"
$ws.Range("C60").Value = "enum"
$ws.Range("D60").Value = "word-paragraph-insert-formatted-text"
$ws.Range("E60").Value = "addPreStyledFormattedText"

# Match the number formatting / vertical alignment used elsewhere in
# the table for this new row.
$ws.Range("A60").NumberFormat = "General"
$ws.Range("B60").NumberFormat = "General"
$ws.Range("D60").NumberFormat = "General"
$ws.Range("E60").NumberFormat = "General"

$ws.Range("B60,D60,E60").VerticalAlignment = -4108

# Restore the view state captured at save time: scrolled so row 33 is
# the first row under the frozen header, with the newly-added
# "getCell" snippet-id cell (now on row 61) selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 33
$ws.Range("D61").Select()
